$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the special fill style that used to be applied to B5:B6 ("Test"
# rows) so every data row goes back to the default (General) style.
$ws.Range("B5:B6").Style = "Normal"

# New data set (dates, apex class / test-class names, category, count)
$data = @(
    @(43808, "AccountTriggerHandler",   "ApexClass", 1),
    @(43809, "PortalApplication",       "ApexClass", 1),
    @(43811, "HttpCreateApprovals",     "ApexClass", 2),
    @(43815, "CaseTriggerHandler",      "ApexClass", 1),
    @(43816, "PortalApplication",       "ApexClass", 2),
    @(43822, "HttpCreateForms",         "ApexClass", 1),
    @(43822, "HttpCreateApprovals",     "ApexClass", 2),
    @(43832, "Test classes",            "ApexClass", 2),
    @(43832, "HttpApprovals",           "ApexClass", 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $entry = $data[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}

# Row 9 is brand new - give its date cell the same date format used by
# the rest of column A.
$ws.Range("A9").NumberFormat = "mm/dd/yy;@"

$ws.Range("G9").Select() | Out-Null
